$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)
$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Lista de sugestões de amizade</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Face</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Sorriso amigável; contato visual constante, mas não excessivamente intenso</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Hostil – Sem sorriso; contato visual intenso</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Indiferente – Sem sorriso; não há muito contato visual</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Voz</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Amigável - caloroso, enérgico</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Hostil - Frio, enérgico</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Indiferente - Frio, menos energético</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Lista de dicas de dominância</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Face</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Dominante – mais contato visual; pode olhar fixamente; desvia o olhar com mais frequência ao ouvir</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Submisso – menos contato visual, mas mais ao ouvir; evita olhar fixamente; mais sorridente</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Corpo</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Dominante - postura aberta e calma; ocupa mais espaço físico; move-se menos; pode usar gestos amplos e empáticos; pode tocar pessoas menos dominantes ocasionalmente</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Submisso - postura nervosa e fechada; pode tocar frequentemente no cabelo, rosto, etc.; pode manter a cabeça baixa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Voz</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Dominante - mais alto; mais controle da conversa</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Submisso - mais suave; segue o exemplo do outro na conversa</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
Write-Host "done"
